# Generate Report for Handoff
# Update status + timestamps on each sheet to reflect handoff completion,
# and widen the Status/date columns that now hold the longer text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-06 18:52:34"

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-06 18:52:29"

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-06 18:52:34"

# --- widen the Status / language columns to fit "Ready for handoff" ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.33   # zh-cn status column
$wsOverview.Columns.Item(6).ColumnWidth = 16.33   # de-de status column
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33        # Status column
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33        # Status column
